# Update the "Förändrad" (Changed) date column C for rows 2-28
# from 45175 (2023-09-06) to 45183 (2023-09-14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45183
    }
}
